# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
# Most cells are plain text (prices/percentages formatted as strings in the source
# sheet), so we write through .Formula (not .Value) to avoid Excel's automatic
# numeric coercion; purely-numeric-looking price strings get a leading apostrophe
# (the standard COM/Excel force-text idiom) so they keep storing as text, matching
# the original inline-string cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "51.821.26"
$ws.Range("E2").Formula = "  -0.69%  "
$ws.Range("D3").Formula = "2.810.58"
$ws.Range("E3").Formula = "  +0.73%  "
$ws.Range("D4").Formula = "'0.999"
$ws.Range("E4").Formula = "  -0.01%  "
$ws.Range("D5").Formula = "'352.59"
$ws.Range("E5").Formula = "  +2.31%  "
$ws.Range("D6").Formula = "'111.16"
$ws.Range("E6").Formula = "  -4.06%  "
$ws.Range("D7").Formula = "'0.561"
$ws.Range("E7").Formula = "  +2.28%  "
$ws.Range("D8").Formula = "'0.999"
$ws.Range("D9").Formula = "'0.596"
$ws.Range("E9").Formula = "  +3.11%  "
$ws.Range("D10").Formula = "'40.56"
$ws.Range("E10").Formula = "  -5.01%  "
$ws.Range("D11").Formula = "'0.0852"
$ws.Range("E12").Formula = "  +0.16%  "
$ws.Range("D13").Formula = "'19.76"
$ws.Range("E13").Formula = "  -1.93%  "
$ws.Range("D14").Formula = "'7.75"
$ws.Range("E14").Formula = "  +0.58%  "
$ws.Range("D15").Formula = "3.243.77"
$ws.Range("E15").Formula = "  +0.69%  "
$ws.Range("D16").Formula = "2.816.85"
$ws.Range("E16").Formula = "  +1.39%  "
$ws.Range("D17").Formula = "'0.916"
$ws.Range("E17").Formula = "  +3.51%  "
$ws.Range("D18").Formula = "51.551.94"
$ws.Range("E18").Formula = "  -0.76%  "
$ws.Range("D19").Formula = "'7.51"
$ws.Range("E19").Formula = "  +6.36%  "
$ws.Range("E20").Formula = "  -4.27%  "
$ws.Range("D21").Formula = "'13.27"
$ws.Range("E21").Formula = "  -0.81%  "
$ws.Range("D22").Formula = "0.0₃0989"
$ws.Range("E22").Formula = "  +0.90%  "
$ws.Range("D23").Formula = "'69.95"
$ws.Range("E23").Formula = "  -0.27%  "
$ws.Range("D24").Formula = "'267.24"
$ws.Range("E24").Formula = "  -1.09%  "
$ws.Range("D25").Formula = "'2.79"
$ws.Range("E25").Formula = "  +0.51%  "
$ws.Range("D26").Formula = "'26.81"
$ws.Range("E26").Formula = "  +0.71%  "
$ws.Range("E27").Formula = "  +0.09%  "
$ws.Range("D28").Formula = "'10.23"
$ws.Range("E28").Formula = "  -0.29%  "
$ws.Range("B29").Formula = "Toncoin"
$ws.Range("C29").Formula = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Formula = "'2.25"
$ws.Range("E29").Formula = "  +0.60%  "
$ws.Range("B30").Formula = "VeChain"
$ws.Range("C30").Formula = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D30").Formula = "'0.0493"
$ws.Range("E30").Formula = "  +20.17%  "
$ws.Range("E31").Formula = "  +0.22%  "
$ws.Range("D32").Formula = "'52.13"
$ws.Range("E32").Formula = "  +3.78%  "
$ws.Range("D33").Formula = "'34.24"
$ws.Range("E33").Formula = "  -1.42%  "
$ws.Range("D34").Formula = "'5.87"
$ws.Range("E34").Formula = "  +2.57%  "
$ws.Range("D35").Formula = "'5.45"
$ws.Range("E35").Formula = "  +9.98%  "
$ws.Range("D36").Formula = "'0.0845"
$ws.Range("E36").Formula = "  +2.78%  "
$ws.Range("D37").Formula = "'0.999"
$ws.Range("E37").Formula = "  +0.04%  "
$ws.Range("E38").Formula = "  +0.78%  "
$ws.Range("D39").Formula = "'2.02"
$ws.Range("E39").Formula = "  -3.82%  "
$ws.Range("D40").Formula = "'18.19"
$ws.Range("E40").Formula = "  -4.32%  "
$ws.Range("D41").Formula = "'0.117"
$ws.Range("D42").Formula = "'126.95"
$ws.Range("E42").Formula = "  -0.04%  "
$ws.Range("D43").Formula = "'23.16"
$ws.Range("E43").Formula = "  -1.09%  "
$ws.Range("D44").Formula = "'2.47"
$ws.Range("E44").Formula = "  -8.14%  "
$ws.Range("E45").Formula = "  -2.25%  "
$ws.Range("D46").Formula = "2.086.75"
$ws.Range("E46").Formula = "  +0.79%  "
$ws.Range("D47").Formula = "'3.31"
$ws.Range("E47").Formula = "  -0.93%  "
$ws.Range("D49").Formula = "'5.94"
$ws.Range("E49").Formula = "  +6.82%  "
$ws.Range("D50").Formula = "'0.975"
$ws.Range("E50").Formula = "  +7.84%  "
$ws.Range("D51").Formula = "'9.01"
$ws.Range("E51").Formula = "  +1.29%  "
